$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the season-record columns (Wins / Losses / Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting (bold, centered, bordered) already used by the other
# header cells by copying the format from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in every player's season record (86-76-0) for rows 2 through 55.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
